# Apply updated "想去人数" (F column) counts to the "展览" and "全部类型" sheets
# as generated by the gh-pages data refresh (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> row number -> new F value
$updates = @{
    "展览" = @{
        2  = 81
        3  = 414
        4  = 3016
        5  = 461
        6  = 37
        9  = 6
        10 = 14382
        13 = 5757
        23 = 2929
        25 = 10540
        27 = 57
        28 = 79
        31 = 66
    }
    "全部类型" = @{
        2  = 81
        3  = 414
        5  = 3016
        6  = 461
        7  = 37
        10 = 6
        11 = 14382
        14 = 5757
        24 = 2929
        27 = 10540
        29 = 57
        30 = 79
        33 = 66
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
